$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 22: "Begynd UC01" ------------------------------------------------
$ws.Range("A22").Value = "Begynd UC01"
$ws.Range("C22").Value = 43965
$ws.Range("D22").Value = 0.41666666666666669
$ws.Range("E22").Value = 0.44791666666666669

# --- Row 23: "Snak med Anders" --------------------------------------------
$ws.Range("A23").Value = "Snak med Anders"
$ws.Range("C23").Value = 43965
$ws.Range("D23").Value = 0.44791666666666669
$ws.Range("E23").Value = 0.5

# --- Row 24: "Pause " -------------------------------------------------------
$ws.Range("A24").Value = "Pause "
$ws.Range("C24").Value = 43965
$ws.Range("D24").Value = 0.5
$ws.Range("E24").Value = 0.54166666666666663

# --- Move the active selection from D22 to A25 -----------------------------
$ws.Range("A25").Select()
